# Update the metadata for indicator 4.2.1 (Кыргызская Республика SDG metadata sheet):
#  - Indicator wording narrowed from "до пяти лет" to "от 36 до 59 месяцев"
#  - Responsible unit renamed from "Отдел статистики домашних хозяйств" to
#    "Управление статистики домашних хозяйств"
#  - New contact person / e-mail / phone / website for the organization

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "4.2.1. Доля детей в возрасте от 36 до 59 месяцев, развивающихся без отклонений в плане здоровья, обучения и психосоциального благополучия, в разбивке по полу"

$ws.Range("B6").Value = "Национальный статистический комитет Кыргызской Республики (Управление статистики домашних хозяйств)"
$ws.Range("B7").Value = "Калымбетова Ы.И."
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

$ws.Range("B6").Select()
